# Insert a new row at position 122 (shifting existing rows 122-185 down to 123-186)
# and populate it with the new week's data, matching the pattern of the
# surrounding rows (same market / category / quality / origin metadata,
# with updated date and price figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(122).Insert()

$ws.Cells.Item(122, 1).Value2  = 9
$ws.Cells.Item(122, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(122, 3).Value2  = "Metropolitana"
$ws.Cells.Item(122, 4).Value2  = 44518
$ws.Cells.Item(122, 5).Value2  = 13
$ws.Cells.Item(122, 6).Value2  = 300000001
$ws.Cells.Item(122, 7).Value2  = "Rabanito"
$ws.Cells.Item(122, 8).Value2  = "Sin especificar"
$ws.Cells.Item(122, 9).Value2  = "Primera"
$ws.Cells.Item(122, 10).Value2 = 8800
$ws.Cells.Item(122, 11).Value2 = 2500
$ws.Cells.Item(122, 12).Value2 = 3000
$ws.Cells.Item(122, 13).Value2 = 2750
$ws.Cells.Item(122, 14).Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(122, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(122, 16).Value2 = 28
$ws.Cells.Item(122, 17).Value2 = 100
$ws.Cells.Item(122, 18).Value2 = "Hortaliza"
